$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the blank separator row (row 2) into the new column I,
# reusing the same style as H2.
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

# Add the new "2021" year header in I3, cloning H3's format (bold
# Times New Roman w/ bottom border) but bumped to size 11.
$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Value = 2021
$ws.Range("I3").Font.Size = 11

# Add the new data value in I4, cloning H4's format (plain Times New
# Roman, no border) but bumped to size 11.
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 149
$ws.Range("I4").Font.Size = 11

# Add the new data value in I5, cloning H5's format (plain Times New
# Roman w/ bottom border) but bumped to size 11.
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = 159
$ws.Range("I5").Font.Size = 11

# Match the saved selection state from the authored workbook.
$ws.Range("K4").Select()
